# Update "three-digit number × one-digit number" problems with new values.
$d = $word.ActiveDocument

$replacements = @(
    @{old="377×4="; new="793×2="},
    @{old="754×8="; new="190×7="},
    @{old="830×2="; new="692×7="},
    @{old="429×3="; new="450×3="},
    @{old="905×7="; new="905×2="},
    @{old="151×9="; new="595×2="},
    @{old="223×9="; new="473×7="},
    @{old="222×4="; new="442×5="},
    @{old="951×5="; new="338×3="},
    @{old="188×7="; new="471×2="},
    @{old="726×2="; new="767×7="},
    @{old="381×5="; new="672×2="},
    @{old="228×9="; new="186×9="},
    @{old="735×9="; new="663×3="},
    @{old="309×3="; new="266×8="},
    @{old="304×2="; new="591×8="},
    @{old="423×8="; new="862×3="},
    @{old="838×7="; new="378×6="},
    @{old="832×5="; new="921×6="},
    @{old="506×4="; new="290×6="},
    @{old="434×7="; new="512×7="},
    @{old="488×2="; new="138×7="},
    @{old="356×7="; new="850×6="},
    @{old="552×5="; new="313×2="},
    @{old="894×4="; new="986×2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
